# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (F column) and occasionally "最低票价" (G column)
# counts across the 展览 / 演出 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 6).Value = 5184
$ws1.Cells.Item(5, 6).Value = 5184
$ws1.Cells.Item(7, 6).Value = 161
$ws1.Cells.Item(8, 6).Value = 212
$ws1.Cells.Item(12, 6).Value = 8651
$ws1.Cells.Item(13, 6).Value = 8651
$ws1.Cells.Item(16, 6).Value = 630
$ws1.Cells.Item(17, 6).Value = 6
$ws1.Cells.Item(18, 6).Value = 2575
$ws1.Cells.Item(20, 6).Value = 2326
$ws1.Cells.Item(23, 6).Value = 2537
$ws1.Cells.Item(25, 6).Value = 18
$ws1.Cells.Item(26, 6).Value = 6508
$ws1.Cells.Item(27, 6).Value = 205
$ws1.Cells.Item(28, 6).Value = 77
$ws1.Cells.Item(31, 6).Value = 463
$ws1.Cells.Item(32, 6).Value = 7052
$ws1.Cells.Item(37, 6).Value = 13
$ws1.Cells.Item(41, 6).Value = 40
$ws1.Cells.Item(45, 6).Value = 72
$ws1.Cells.Item(47, 6).Value = 69
$ws1.Cells.Item(48, 6).Value = 540
$ws1.Cells.Item(49, 6).Value = 3116
$ws1.Cells.Item(51, 6).Value = 1127

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 13
$ws2.Cells.Item(4, 6).Value = 22
$ws2.Cells.Item(5, 6).Value = 191
$ws2.Cells.Item(5, 7).Value = 180
$ws2.Cells.Item(7, 6).Value = 83
$ws2.Cells.Item(14, 6).Value = 159

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 5184
$ws4.Cells.Item(4, 6).Value = 5184
$ws4.Cells.Item(6, 6).Value = 161
$ws4.Cells.Item(7, 6).Value = 212
$ws4.Cells.Item(11, 6).Value = 8651
$ws4.Cells.Item(12, 6).Value = 8651
$ws4.Cells.Item(15, 6).Value = 630
$ws4.Cells.Item(16, 6).Value = 2575
$ws4.Cells.Item(17, 6).Value = 22
$ws4.Cells.Item(18, 6).Value = 191
$ws4.Cells.Item(18, 7).Value = 180
$ws4.Cells.Item(20, 6).Value = 2326
$ws4.Cells.Item(21, 6).Value = 83
$ws4.Cells.Item(22, 6).Value = 2537
$ws4.Cells.Item(25, 6).Value = 18
$ws4.Cells.Item(26, 6).Value = 6508
$ws4.Cells.Item(27, 6).Value = 205
$ws4.Cells.Item(29, 6).Value = 77
$ws4.Cells.Item(32, 6).Value = 463
$ws4.Cells.Item(33, 6).Value = 7052
$ws4.Cells.Item(36, 6).Value = 13
$ws4.Cells.Item(38, 6).Value = 40
$ws4.Cells.Item(44, 6).Value = 69
$ws4.Cells.Item(45, 6).Value = 540
$ws4.Cells.Item(46, 6).Value = 159
$ws4.Cells.Item(47, 6).Value = 3116
$ws4.Cells.Item(50, 6).Value = 1127

$wb.Save()
